# Apply the "updated sheet as per main" commit:
#  - Insert a new "ProductForOrder" worksheet right after "MembershipStatus"
#    with ProductCategory/ProductName header + one Kanji Kentei product row.
#  - Update the "NewUser" sample credentials to the later (17082021184752) set.
#  - Extend "MembershipStatus" with a new header row + CASEC/Kanji Kentei columns.

$wb = $excel.ActiveWorkbook

# --- 1. New "ProductForOrder" sheet, positioned after "MembershipStatus" ---
$membership = $wb.Worksheets.Item("MembershipStatus")
$productForOrder = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $membership)
$productForOrder.Name = "ProductForOrder"

$productForOrder.Range("A1").Value = "ProductCategory"
$productForOrder.Range("B1").Value = "ProductName"
$productForOrder.Range("A2").Value = "漢検"
$productForOrder.Range("B2").Value = "スタギア漢検３級プレミアム６ヶ月"

# --- 2. "NewUser" sheet: refresh the sample account to the newer test user ---
$newUser = $wb.Worksheets.Item("NewUser")
$newUser.Range("A2").Value = "TestPF1221+17082021184752@gmail.com"
$newUser.Range("C2").Value = "TestPF1221_17082021184752"

# --- 3. "MembershipStatus" sheet: push the existing header row down and add
#        a new header row above it plus two new CASEC / Kanji Kentei columns ---
$membershipStatus = $wb.Worksheets.Item("MembershipStatus")
$membershipStatus.Range("A2").Value = $membershipStatus.Range("A1").Value()
$membershipStatus.Range("B2").Value = $membershipStatus.Range("B1").Value()

$membershipStatus.Range("A1").Value = "Freeuser"
$membershipStatus.Range("B1").Value = "Primeuser"
$membershipStatus.Range("C1").Value = "CASEC"
$membershipStatus.Range("D1").Value = "Kanji Kentei"
$membershipStatus.Range("C2").Value = "CASEC"
$membershipStatus.Range("D2").Value = "漢検"
$membershipStatus.Range("K21").Select()

# --- 4. Make the new sheet the active tab, matching the authored workbook view ---
$productForOrder.Activate()
$productForOrder.Range("J14").Select()
